$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential disclosure text date (2021-03-23 -> 2021-03-24)
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."

# Update weight (D) and percent change (E) values for rows 2-13
$ws.Range("D2").Value = 0.03030213330585698
$ws.Range("E2").Value = -0.004069552349241556

$ws.Range("D3").Value = 0.02400041516846191
$ws.Range("E3").Value = 0.0007400098667982302

$ws.Range("D4").Value = 0.05244575287449998
$ws.Range("E4").Value = -0.00489955903968653

$ws.Range("D5").Value = 0.1382692362437029
$ws.Range("E5").Value = -0.007231404958677579

$ws.Range("D6").Value = 0.03038564565119854
$ws.Range("E6").Value = 0.0248062015503876

$ws.Range("D7").Value = 0.1204100645098228
$ws.Range("E7").Value = -0.01820114255347394

$ws.Range("D8").Value = 0.09956913675946405
$ws.Range("E8").Value = 0.004997001798920708

$ws.Range("D9").Value = 0.0272764167938666
$ws.Range("E9").Value = 0.00612341026848795

$ws.Range("D10").Value = 0.1201670851521096
$ws.Range("E10").Value = 0.003144654088050425

$ws.Range("D11").Value = 0.2507423529521425
$ws.Range("E11").Value = -0.01498520568865125

$ws.Range("D12").Value = 0.1064317605888741
$ws.Range("E12").Value = -0.02296450939457195

$ws.Range("E13").Value = -0.007959368913706966

# Restore worksheet protection that was temporarily lifted to allow the edits above
$ws.Protect()
